# Replace the red-highlighted "?????" placeholder at the top of the
# document with the finalized date "30th June 2021" (no highlight),
# splitting "th" into its own superscript run.

$d = $word.ActiveDocument

# Locate the placeholder run.
$rng = $d.Content
$ok = $rng.Find.Execute("?????")

if ($ok) {
    # Replace the placeholder text with "30" (keeps this run's existing
    # bold/underline/size formatting, including the complex-script
    # counterparts) and then strip the red highlight from just this run.
    $rng.Text = "30"
    $rng.HighlightColorIndex = 0

    # Insert a new run right after "30" for the ordinal suffix "th".
    $thStart = $rng.End
    $insAfter30 = $d.Range($thStart, $thStart)
    $insAfter30.InsertAfter("th")

    # Give the new "th" run the same character formatting as "30" by
    # cloning its FormattedText, then fix the text back to "th" and mark
    # it superscript.
    $thRange = $d.Range($thStart, $thStart + 2)
    $thRange.FormattedText = $rng.FormattedText
    $thRange = $d.Range($thStart, $thStart + 2)
    $thRange.Text = "th"
    $thRange.Font.Superscript = $true

    # Insert a final run for " June 2021" using the same base formatting
    # as "30" (no superscript, no highlight).
    $juneStart = $thRange.End
    $insAfterTh = $d.Range($juneStart, $juneStart)
    $insAfterTh.InsertAfter(" June 2021")

    $juneLen = (" June 2021").Length
    $juneRange = $d.Range($juneStart, $juneStart + $juneLen)
    $juneRange.FormattedText = $rng.FormattedText
    # $rng.FormattedText is only 2 characters ("30"), so after the
    # assignment above the range collapses down to those 2 characters;
    # re-fetch it at that (shorter) width before restoring the real text.
    $juneRange = $d.Range($juneStart, $juneStart + 2)
    $juneRange.Text = " June 2021"
}
